$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("F2").Value = 2.48
$ws.Range("G2").Value = 2.62
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 3.4

# Row 3 updates
$ws.Range("G3").Value = 1.95
$ws.Range("I3").Value = 6.6
$ws.Range("K3").Value = 3.5
$ws.Range("P3").Value = 1.5
$ws.Range("Q3").Value = 2.64
